$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refit results: a handful of D/E (Prediction/Error) values shifted slightly ---
$ws.Range("D2").Value = 0.5742796834542762
$ws.Range("E2").Value = 0.5742796834542762

$ws.Range("D4").Value = 0.5678397472716558
$ws.Range("E4").Value = 0.5678397472716558

$ws.Range("D6").Value = 0.5756907986063734
$ws.Range("E6").Value = 0.5756907986063734

$ws.Range("D7").Value = 0.5723472412718229
$ws.Range("E7").Value = 0.4276527587281771

$ws.Range("D11").Value = 0.5559417284434623
$ws.Range("E11").Value = 0.4440582715565377

# --- New "Label" column (H): 0 = Control patient, 1 = MDD patient ---
# Give H1 the same header formatting (bold / centered / bordered) as the
# other header cells before writing its text.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Label"

$labels = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1;
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1
}

foreach ($row in 2..21) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}
